# STOR_QTR_FIN.xlsx -- "Doing Updates for Financials"
# Two new quarterly columns (31-Dec-18 and 30-Sep-18) are inserted at the
# front of each of the three statements (Income Statement, Balance Sheet,
# Cash Flow Statement) on the STOR sheet. Inserting the columns pushes the
# existing quarters (old D:K) right to F:M, and the two freshly opened
# columns (D:E) are then populated with the newly reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STOR")

# Insert two blank columns at D:E -- everything that used to live in D:K
# shifts right to F:M, formulas/format intact.
$ws.Range("D:E").EntireColumn.Insert()

# The freshly inserted columns don't inherit number formatting from either
# neighbour automatically, so pull it from column F (the data that used to
# be column D before the shift) across the whole data block.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# row, new column D value, new column E value
$newData = @(
    @(7, 43465, 43373),
    @(8, 146700, 137000),
    @(9, 1400, 800),
    @(10, 145300, 136200),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 5200, "NA"),
    @(15, 49500, 45800),
    @(17, 68600, 58000),
    @(18, 78100, 79000),
    @(20, 0, 0),
    @(21, 127600, 124700),
    @(22, 36000, 31800),
    @(23, 42100, 47100),
    @(24, 300, 100),
    @(25, 0, 0),
    @(26, 41800, 47000),
    @(27, 41700, 46900),
    @(28, 0, 0),
    @(29, 14800, 1200),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 0, 0),
    @(33, 56500, 48100),
    @(34, 0, 0),
    @(35, 56500, 48100),
    @(38, 43465, 43373),
    @(41, 27500, 25600),
    @(42, 0, 0),
    @(43, 0, 0),
    @(44, 0, 0),
    @(45, 0, 0),
    @(46, 0, 0),
    @(47, 351200, 352000),
    @(48, 6668000, 6311600),
    @(49, 0, 0),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 0, 0),
    @(53, 0, 0),
    @(54, 7114000, 6751700),
    @(57, 117200, 106700),
    @(58, 0, 0),
    @(59, 73000, 69900),
    @(60, 0, 0),
    @(61, 3060300, 2956400),
    @(62, 0, 0),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 3250500, 3133000),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, -267700, -250200),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 3863500, 3618700),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 56500, 48100),
    @(83, 49500, 45800),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 108600, 103600),
    @(91, -461000, -470700),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -397200, -488400),
    @(96, -69900, -63600),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, 298600, 361200),
    @(101, 0, 0),
    @(102, 10100, -23600)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value = $entry[1]
    $ws.Cells.Item($r, 5).Value = $entry[2]
}
